$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.639.26"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "2.292.60"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "96.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "268.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.617"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.76%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.612"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.80"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0936"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.90"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.105"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.29%  "
$ws.Range("D14").Value = "2.636.00"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("D17").Value = "2.292.06"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").Value = "43.606.98"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("E19").Value = "  +2.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +12.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.05%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.90%  "
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("E28").Value = "  +2.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0899"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.73%  "
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("E36").Value = "  -2.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0354"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.78%  "
$ws.Range("E38").Value = "  -3.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.242"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.31"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.16%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "65.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.08%  "
$ws.Range("E46").Value = "  -4.04%  "
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "97.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.79%  "
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.186"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.04%  "
